$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '68.356.24'
$ws.Cells.Item(2, 5).Value = '  +0.48%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.637.11'
$ws.Cells.Item(3, 5).Value = '  -0.15%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.03%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''600.23'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +0.60%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''154.41'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +1.24%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.04%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.546'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +0.24%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.636.55'
$ws.Cells.Item(9, 5).Value = '  -0.15%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''0.145'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +6.93%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.54%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +0.63%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''0.352'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +0.93%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''28.11'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +1.57%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''0.0000192'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -0.25%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '3.120.90'
$ws.Cells.Item(16, 5).Value = '  -0.26%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '68.202.28'
$ws.Cells.Item(17, 5).Value = '  +0.33%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.630.63'
$ws.Cells.Item(18, 5).Value = '  -0.42%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''11.50'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +2.58%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''365.65'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -2.72%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''7.50'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +0.38%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''4.38'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +3.42%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''4.87'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +0.83%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''2.08'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +1.40%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''73.62'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -0.64%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''0.999'
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -0.14%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''9.87'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -1.21%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'PEPE'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(28, 4).Value = '''0.0000106'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +1.17%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'WrappedeETH'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(29, 4).Value = '2.772.06'
$ws.Cells.Item(29, 5).Value = '  -0.18%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''0.997'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -0.37%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''574.38'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -2.91%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(32, 4).Value = '''8.07'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +3.54%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'Fetch.AI'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(33, 4).Value = '''1.43'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +3.35%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +1.35%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''0.132'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +4.00%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -0.08%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''1.59'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +4.26%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''160.34'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +1.17%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''19.44'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +0.86%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''1.89'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -0.40%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +0.50%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''5.43'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +2.19%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +1.34%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''17.72'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +3.47%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +9.66%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +0.01%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''40.58'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +0.53%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''157.38'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +0.79%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '''3.78'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +2.35%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''1.71'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +0.56%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''21.93'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +2.80%  '
